$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Make room for the new "Indeks" column. This shifts the existing
#    "UL" / "Last delivery" columns from A/B to D/E (old data + old styles
#    travel with them), and leaves A:C genuinely blank/unformatted.
# ---------------------------------------------------------------------------
$ws.Columns("A:C").Insert()

# A completely untouched scratch cell we use as a "format eraser": pasting
# its (default) format onto another cell resets that cell back to the
# workbook's implicit default style (no explicit <s> attribute at all).
$blank = $ws.Range("ZZ1")

# ---------------------------------------------------------------------------
# 2. Header row: Indeks / UL / Last delivery
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Indeks"
$ws.Range("B1").Value = "UL"
$ws.Range("C1").Value = "Last delivery"

# ---------------------------------------------------------------------------
# 3. Data rows (new sort order: by date, descending -- already reflected by
#    the literal values below) plus the new row 6.
# ---------------------------------------------------------------------------
$indeks = @(37221, 56138, 55472, 59958, 60988)
$ul     = @(473221, 345919, 132781, 174211, 360486)
$dates  = @("2024-08-20", "2024-12-03", "2025-03-03", "2025-03-13", "2025-05-20")

for ($i = 0; $i -lt 5; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $indeks[$i]
    $ws.Cells.Item($r, 2).Value = $ul[$i]

    # Force the date-like text to be stored as plain text (matching the
    # original file, which keeps these as inline strings, not real dates).
    $dcell = $ws.Cells.Item($r, 3)
    $dcell.NumberFormat = "@"
    $dcell.Value = $dates[$i]
    # Wipe the leftover "@" number format back to the sheet default so the
    # cell ends up with no explicit style, same as its neighbours.
    $blank.Copy()
    $dcell.PasteSpecial(-4122)
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------------
# 4. Header style: bold, centered, bordered (same look as the pre-existing,
#    already-defined-but-unused header format). Derive it from the OLD
#    header cell (now shifted to D1, still carrying that exact font/border)
#    so we reuse the workbook's existing style entry instead of minting a
#    brand new one.
# ---------------------------------------------------------------------------
$ws.Range("D1").Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A1:C1").WrapText = $false
$ws.Range("A1:C1").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 5. Data cells should carry no special formatting (no wrap, no border) --
#    reset them to the implicit default style via the same blank-paste trick.
# ---------------------------------------------------------------------------
$blank.Copy()
$ws.Range("A2:B6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 6. Drop the old (now empty of meaning) UL / Last delivery columns.
# ---------------------------------------------------------------------------
$ws.Columns("D:E").Delete()

# ---------------------------------------------------------------------------
# 7. Clear the custom row heights inherited from the original sheet.
# ---------------------------------------------------------------------------
$ws.Range("A1:C6").EntireRow.AutoFit()

$blank.Clear()
